$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells keep their text representation (avoid numeric
# auto-conversion stripping trailing zeros / using scientific notation).
$ws.Range('B2:E51').NumberFormat = '@'

$ws.Range('D2').Value = '24.726.77'
$ws.Range('E2').Value = '  +1.50%  '
$ws.Range('D3').Value = '1.698.62'
$ws.Range('E3').Value = '  +1.46%  '
$ws.Range('D4').Value = '1.008'
$ws.Range('E4').Value = '  +0.29%  '
$ws.Range('D5').Value = '311.10'
$ws.Range('E5').Value = '  +1.87%  '
$ws.Range('D7').Value = '0.3719'
$ws.Range('E7').Value = '  +1.00%  '
$ws.Range('D8').Value = '49.19'
$ws.Range('E8').Value = '  +3.05%  '
$ws.Range('D9').Value = '0.3414'
$ws.Range('E9').Value = '  -0.28%  '
$ws.Range('D10').Value = '1.209'
$ws.Range('E10').Value = '  +4.40%  '
$ws.Range('D11').Value = '0.07450'
$ws.Range('E11').Value = '  +3.47%  '
$ws.Range('D12').Value = '1.004'
$ws.Range('E12').Value = '  +0.18%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = '6.307'
$ws.Range('E13').Value = '  +2.88%  '
$ws.Range('B14').Value = 'Solana'
$ws.Range('C14').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D14').Value = '20.87'
$ws.Range('E14').Value = '  +3.91%  '
$ws.Range('D15').Value = '6.983'
$ws.Range('E15').Value = '  +3.83%  '
$ws.Range('D16').Value = '1.699.84'
$ws.Range('E16').Value = '  +1.48%  '
$ws.Range('D17').Value = '0.00001121'
$ws.Range('E17').Value = '  +1.95%  '
$ws.Range('D18').Value = '0.06693'
$ws.Range('E18').Value = '  +0.58%  '
$ws.Range('E19').Value = '  +0.23%  '
$ws.Range('D20').Value = '83.21'
$ws.Range('E20').Value = '  +3.48%  '
$ws.Range('D21').Value = '17.10'
$ws.Range('E21').Value = '  +4.04%  '
$ws.Range('D22').Value = '6.314'
$ws.Range('E22').Value = '  +3.62%  '
$ws.Range('D23').Value = '12.92'
$ws.Range('E23').Value = '  +6.38%  '
$ws.Range('D24').Value = '24.755.02'
$ws.Range('E24').Value = '  +1.80%  '
$ws.Range('D25').Value = '2.455'
$ws.Range('E25').Value = '  +0.78%  '
$ws.Range('D26').Value = '2.758'
$ws.Range('E26').Value = '  +4.23%  '
$ws.Range('D27').Value = '20.18'
$ws.Range('E27').Value = '  +4.08%  '
$ws.Range('D28').Value = '148.75'
$ws.Range('E28').Value = '  -2.76%  '
$ws.Range('B29').Value = 'ImmutableX'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D29').Value = '1.258'
$ws.Range('E29').Value = '  +29.75%  '
$ws.Range('B30').Value = 'BitcoinCash'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D30').Value = '131.50'
$ws.Range('E30').Value = '  +3.18%  '
$ws.Range('B31').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C31').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D31').Value = '1.886.41'
$ws.Range('E31').Value = '  +1.41%  '
$ws.Range('D32').Value = '6.713'
$ws.Range('E32').Value = '  +7.41%  '
$ws.Range('D33').Value = '4.213'
$ws.Range('E33').Value = '  +4.47%  '
$ws.Range('B34').Value = 'WEMIXTOKEN'
$ws.Range('C34').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D34').Value = '1.758'
$ws.Range('E34').Value = '  +4.45%  '
$ws.Range('B35').Value = 'Aptos'
$ws.Range('C35').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D35').Value = '13.53'
$ws.Range('E35').Value = '  +9.75%  '
$ws.Range('B36').Value = 'Stellar'
$ws.Range('C36').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D36').Value = '0.08685'
$ws.Range('E36').Value = '  +3.05%  '
$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D37').Value = '0.06621'
$ws.Range('E37').Value = '  +3.85%  '
$ws.Range('B38').Value = 'InternetComputer(DFINITY)'
$ws.Range('C38').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D38').Value = '5.530'
$ws.Range('E38').Value = '  +4.24%  '
$ws.Range('D39').Value = '9.006'
$ws.Range('E39').Value = '  +3.77%  '
$ws.Range('D40').Value = '0.02393'
$ws.Range('E40').Value = '  +3.80%  '
$ws.Range('E41').Value = '  +6.15%  '
$ws.Range('D42').Value = '1.267'
$ws.Range('E42').Value = '  +2.08%  '
$ws.Range('D43').Value = '0.6380'
$ws.Range('E43').Value = '  +5.02%  '
$ws.Range('E44').Value = '  +0.25%  '
$ws.Range('D45').Value = '13.79'
$ws.Range('E45').Value = '  +6.68%  '
$ws.Range('D46').Value = '0.6080'
$ws.Range('E46').Value = '  +3.63%  '
$ws.Range('D47').Value = '3.813'
$ws.Range('E47').Value = '  +1.60%  '
$ws.Range('D48').Value = '2.109'
$ws.Range('E48').Value = '  +4.74%  '
$ws.Range('D49').Value = '128.92'
$ws.Range('E49').Value = '  +2.78%  '
$ws.Range('D50').Value = '0.07250'
$ws.Range('E50').Value = '  +1.52%  '
$ws.Range('D51').Value = '79.31'
$ws.Range('E51').Value = '  +4.88%  '
